# Auto-generated Excel COM-interop script
# Applies numeric "want to go" count updates (F column) across all four sheets,
# and appends one new row of data to "本地生活" (and reflects it in "全部类型"'s pre-existing row).

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 136
$ws.Range("F4").Value = 909
$ws.Range("F5").Value = 1077
$ws.Range("F7").Value = 333
$ws.Range("F8").Value = 677
$ws.Range("F9").Value = 12224
$ws.Range("F11").Value = 2177
$ws.Range("F13").Value = 257
$ws.Range("F15").Value = 1231
$ws.Range("F16").Value = 214
$ws.Range("F20").Value = 301
$ws.Range("F21").Value = 2924
$ws.Range("F22").Value = 761
$ws.Range("F23").Value = 4152
$ws.Range("F24").Value = 4152
$ws.Range("F25").Value = 1123
$ws.Range("F26").Value = 864
$ws.Range("F30").Value = 1057
$ws.Range("F31").Value = 53
$ws.Range("F32").Value = 103
$ws.Range("F36").Value = 29
$ws.Range("F38").Value = 4430
$ws.Range("F39").Value = 16
$ws.Range("F40").Value = 4554
$ws.Range("F41").Value = 5545
$ws.Range("F44").Value = 70
$ws.Range("F46").Value = 328
$ws.Range("F47").Value = 81
$ws.Range("F49").Value = 4113
$ws.Range("F50").Value = 132

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 7
$ws.Range("F13").Value = 1038
$ws.Range("F22").Value = 14

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 461
$ws.Range("F4").Value = 81

# Append new row 5 (new 漫展/活动 entry)
$ws.Range("A5").Value = 4
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = "'2024-10-06"
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").Value = "杭州· LoveLive! Series Asia Tour 2024~伴你圆梦~直播"
$ws.Range("D5").Value = "通货路918号粮仓艺术公园7号楼 SoFunLivehouse"
$ws.Range("E5").Value = "2024.10.06 18:45-10.06 21:30"
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 250
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92903"
$ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202409/dvcSAoFl1727183785196.jpeg"

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 461
$ws.Range("F4").Value = 81
$ws.Range("F5").Value = 136
$ws.Range("F6").Value = 909
$ws.Range("F8").Value = 333
$ws.Range("F9").Value = 677
$ws.Range("F10").Value = 12224
$ws.Range("F11").Value = 2177
$ws.Range("F12").Value = 257
$ws.Range("F13").Value = 1231
$ws.Range("F16").Value = 301
$ws.Range("F17").Value = 2924
$ws.Range("F18").Value = 761
$ws.Range("F20").Value = 4152
$ws.Range("F21").Value = 1123
$ws.Range("F23").Value = 864
$ws.Range("F27").Value = 1057
$ws.Range("F28").Value = 53
$ws.Range("F29").Value = 103
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 4430
$ws.Range("F34").Value = 4554
$ws.Range("F39").Value = 328
$ws.Range("F42").Value = 81
$ws.Range("F44").Value = 4113
$ws.Range("F47").Value = 14
$ws.Range("F49").Value = 132

